# Update the dSF column (F) values as part of a data re-pull / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -3
    4  = -2
    5  = 1
    6  = -8
    10 = -8
    11 = 0
    12 = -2
    13 = -3
    14 = 2
    15 = 2
    16 = 4
    17 = 3
    18 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
